$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.ClearFormats()
}

Set-TextValue "D2" '59.409.10'
Set-TextValue "E2" '  +0.08%  '

Set-TextValue "D3" '2.607.85'
Set-TextValue "E3" '  -0.65%  '

Set-TextValue "E4" '  -0.26%  '

Set-TextValue "D5" '541.39'
Set-TextValue "E5" '  +3.78%  '

Set-TextValue "D6" '141.24'
Set-TextValue "E6" '  +0.06%  '

Set-TextValue "E7" '  +0.13%  '

Set-TextValue "D8" '0.568'
Set-TextValue "E8" '  +0.19%  '

Set-TextValue "D9" '6.46'
Set-TextValue "E9" '  -1.53%  '

Set-TextValue "E10" '  +1.74%  '

Set-TextValue "E11" '  +0.87%  '

Set-TextValue "D13" '3.067.36'
Set-TextValue "E13" '  -0.43%  '

Set-TextValue "D14" '59.331.97'
Set-TextValue "E14" '  +0.04%  '

Set-TextValue "D15" '20.59'
Set-TextValue "E15" '  +0.68%  '

Set-TextValue "D16" '2.651.54'
Set-TextValue "E16" '  +1.61%  '

Set-TextValue "D17" '0.0000134'
Set-TextValue "E17" '  +0.60%  '

Set-TextValue "D18" '343.54'
Set-TextValue "E18" '  +1.34%  '

Set-TextValue "D19" '4.35'
Set-TextValue "E19" '  +0.46%  '

Set-TextValue "D20" '10.13'
Set-TextValue "E20" '  -0.81%  '

Set-TextValue "D21" '6.43'
Set-TextValue "E21" '  -1.91%  '

Set-TextValue "E22" '  +0.22%  '

Set-TextValue "D23" '67.35'
Set-TextValue "E23" '  +1.32%  '

Set-TextValue "E24" '  -0.89%  '

Set-TextValue "D25" '0.408'
Set-TextValue "E25" '  +0.73%  '

Set-TextValue "E26" '  +0.20%  '

Set-TextValue "D27" '7.20'
Set-TextValue "E27" '  +1.80%  '

Set-TextValue "E28" '  +0.14%  '

Set-TextValue "D29" '0.0₃0739'
Set-TextValue "E29" '  +1.80%  '

Set-TextValue "D30" '1.67'
Set-TextValue "E30" '  +6.52%  '

Set-TextValue "D31" '5.84'
Set-TextValue "E31" '  -1.95%  '

Set-TextValue "D32" '18.78'
Set-TextValue "E32" '  -0.13%  '

Set-TextValue "D33" '149.41'
Set-TextValue "E33" '  +0.16%  '

Set-TextValue "D34" '3.98'
Set-TextValue "E34" '  -0.74%  '

Set-TextValue "E35" '  -1.34%  '

Set-TextValue "D36" '37.03'
Set-TextValue "E36" '  +1.86%  '

Set-TextValue "E37" '  -0.27%  '

Set-TextValue "D38" '0.834'
Set-TextValue "E38" '  +0.21%  '

Set-TextValue "D39" '0.813'
Set-TextValue "E39" '  -1.20%  '

Set-TextValue "E40" '  +0.34%  '

Set-TextValue "E41" '  +0.17%  '

Set-TextValue "D42" '274.41'
Set-TextValue "E42" '  -0.91%  '

$ws.Range("B43").Value = 'Mantle'
$ws.Range("C43").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue "D43" '0.597'
Set-TextValue "E43" '  +0.80%  '

$ws.Range("B44").Value = 'WhiteBITCoin'
$ws.Range("C44").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
Set-TextValue "D44" '10.76'
Set-TextValue "E44" '  +0.19%  '

Set-TextValue "D45" '0.0957'
Set-TextValue "E45" '  +0.33%  '

Set-TextValue "D46" '0.0524'
Set-TextValue "E46" '  +0.51%  '

Set-TextValue "D47" '1.948.08'
Set-TextValue "E47" '  -2.01%  '

Set-TextValue "E48" '  +0.96%  '

Set-TextValue "D49" '4.54'
Set-TextValue "E49" '  +0.71%  '

Set-TextValue "D50" '18.30'
Set-TextValue "E50" '  +1.09%  '

Set-TextValue "D51" '111.09'
Set-TextValue "E51" '  -2.35%  '
